$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to write: Row, A, B, C, D, E, F, G, H, I, J, K (K may be $null -> leave blank)
$rows = @(
    @(8, 6, "11-12-2025", "MELVIN", 9496952220, "01-01-2026", "Vishnu Raj", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "KURTHA BULK SIZE NOT CONFIRM VISIT AGAIN"),
    @(9, 7, "11-12-2025", "manu", 8590372779, "28-12-2025", "ABIN KRISHNA S", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", ".."),
    @(10, 8, "11-12-2025", "sayid", 8590825100, "25-01-2026", "ABIN KRISHNA S", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY", "-", ".."),
    @(11, 9, "11-12-2025", "Amith", 8281546625, "31-12-2025", "Althaf A", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "they will inform tomorrow"),
    @(12, 10, "11-12-2025", "afsal", 9895939667, "08-01-2026", "AKSHAYDAS.M", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "-"),
    @(13, 11, "11-12-2025", "Ajumal", 8089744330, "01-01-2026", "SWAROOP S K", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "."),
    @(14, 12, "12-12-2025", "sharath", 8075922593, "14-12-2025", "AKSHAYDAS.M", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", ".."),
    @(15, 13, "12-12-2025", "joseph", 6282396016, "05-01-2026", "ABIN KRISHNA S", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", "48size"),
    @(16, 14, "13-12-2025", "adhnana", 8304923273, "13-12-2025", "Vishnu prasad. P", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "size issue"),
    @(17, 15, "13-12-2025", "vivek", 8891349992, "10-01-2026", "ABIN KRISHNA S", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", ".."),
    @(18, 16, "13-12-2025", "jiby", 9946814700, "17-12-2025", "ABIN KRISHNA S", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", ".."),
    @(19, 17, "13-12-2025", "anwar", 6282111567, "21-12-2025", "ABIN KRISHNA S", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", ".."),
    @(20, 18, "13-12-2025", "cristy", 7907562960, "19-12-2025", "AKSHAYDAS.M", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", ".."),
    @(21, 19, "14-12-2025", "Arun", 9937725569, "31-12-2025", "Althaf A", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", "he didn’t suit the pant size"),
    @(22, 20, "14-12-2025", "Jithin", 9562958782, "27-12-2025", "Althaf A", "Loss", "PRODUCT", "REQUIRED DESIGN NOT AVAILABLE", "-", "bulk order"),
    @(23, 21, "14-12-2025", "Nithin", 9188031437, "26-01-2026", "RAMITH K", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", "Loss"),
    @(24, 22, "14-12-2025", "melvin", 9895578570, "29-12-2025", "ABIN KRISHNA S", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", ".."),
    @(25, 23, "15-12-2025", "Tojin", 8943034163, "10-01-2026", "Althaf A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", $null),
    @(26, 24, "15-12-2025", "anush", 8086462505, "18-01-2026", "AKSHAYDAS.M", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", "kurutha"),
    @(27, 25, "16-12-2025", "Abraham", 9747793966, "10-01-2026", "Althaf A", "Loss", "PRICING", "RENT TO HIGH", "-", "."),
    @(28, 26, "16-12-2025", "Ashik", 9037260665, "22-12-2025", "Althaf A", "Loss", "PRODUCT", "REQUIRED DESIGN NOT AVAILABLE", "-", "."),
    @(29, 27, "16-12-2025", "sandeep", 9947153224, "16-12-2025", "Nithin binoy", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", "loos"),
    @(30, 28, "17-12-2025", "Jebin", 8124004551, "14-01-2026", "Althaf A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "they will come sunday"),
    @(31, 29, "17-12-2025", "benadict", 9562474208, "17-12-2025", "Nithin binoy", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE", "-", "loss"),
    @(32, 30, "17-12-2025", "toy", 9846105212, "28-12-2025", "SWAROOP S K", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", "."),
    @(33, 31, "17-12-2025", "suhail", 8137912784, "17-12-2025", "Mohamed Farseen C S", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "enquiry"),
    @(34, 32, "17-12-2025", "althaf", 7306166226, "17-12-2025", "Nithin binoy", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE", "-", "loss"),
    @(35, 33, "17-12-2025", "naveen", 8281044850, "08-01-2026", "ABIN KRISHNA S", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE", "-", ".."),
    @(36, 34, "18-12-2025", "Umar", 9447767211, "31-12-2025", "Althaf A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "family dissaproval they bought new one instead of rental"),
    @(37, 35, "18-12-2025", "vishnu", 8281533975, "18-12-2025", "Mohamed Farseen C S", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "family disapproved entire model in our store"),
    @(38, 36, "18-12-2025", "abin", 9895568883, "18-12-2025", "Nithin binoy", "Loss", "ENQUIRY", "SHOE", "-", "size not available 10 or 9 black glossy"),
    @(39, 37, "18-12-2025", "Gokul", 8129932346, "11-01-2026", "Althaf A", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "internal issues"),
    @(40, 38, "18-12-2025", "Amal", 8921643216, "18-12-2025", "Althaf A", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", "texture black NP 38 not available another booking"),
    @(41, 39, "18-12-2025", "isimal", 9895447752, "29-12-2025", "SWAROOP S K", "Loss", "PRODUCT", "Product Already Booked", "-", "prem black suit size required 36 set."),
    @(42, 40, "18-12-2025", "vishnu", 9995931223, "18-12-2025", "Mohamed Farseen C S", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", "product not available in our store  already booked requirement for premium black shall"),
    @(43, 41, "18-12-2025", "mini", 9497680722, "18-12-2025", "Mohamed Farseen C S", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", "black texture non premium"),
    @(44, 42, "18-12-2025", "Don jose", 9895574884, "18-12-2025", "Mohamed Farseen C S", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE", "-", "customer required formal kids suite"),
    @(45, 43, "19-12-2025", "suthish", 8606277718, "10-01-2026", "SWAROOP S K", "Loss", "PRODUCT", "Product Already Booked", "-", "blue shiny 42 size blazer"),
    @(46, 44, "19-12-2025", "khaleel", 8089977280, "19-12-2025", "Mohamed Farseen C S", "Loss", "CUSTOMER INTERNAL ISSUES", "BUDGET RESTRICTIONS", "-", "bride restrictions"),
    @(47, 45, "19-12-2025", "Rizwan", 8086945836, "25-01-2026", "Althaf A", "Loss", "PRODUCT", "REQUIRED DESIGN NOT AVAILABLE", "-", "black with minimal print"),
    @(48, 46, "19-12-2025", "Akshay", 9400891463, "09-02-2026", "Althaf A", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "family disapproval"),
    @(49, 47, "20-12-2025", "musfar", 9061752911, "10-01-2026", "SWAROOP S K", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "."),
    @(50, 48, "20-12-2025", "Zayan", 9846858838, "21-12-2025", "Binu Augustine", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE", "-", "function is tommorow"),
    @(51, 49, "20-12-2025", "safeer", 8606123654, "23-12-2025", "SWAROOP S K", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "."),
    @(52, 50, "20-12-2025", "keba paul", 8590118941, "26-12-2025", "ABIN KRISHNA S", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "size not availabe"),
    @(53, 51, "20-12-2025", "linu", 7293168811, "27-12-2025", "ABIN KRISHNA S", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", "enquiry"),
    @(54, 52, "20-12-2025", "indrajith", 7012239991, "06-05-2026", "ABIN KRISHNA S", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", "for enquiry"),
    @(55, 53, "21-12-2025", "Anoop", 7974143534, "04-01-2026", "SWAROOP S K", "Loss", "PRODUCT", "Product Already Booked", "-", "coustomer is conformed product but product is booked. Product - Cream open B/G"),
    @(56, 54, "21-12-2025", "binoy", 9061762422, "21-12-2025", "Mohamed Farseen C S", "Loss", "ENQUIRY", "ENQUIRY WITHOUT TRIAL", "-", "enquiry for open bengala customer not satisfied"),
    @(57, 55, "21-12-2025", "pranav", 7012395349, "21-12-2025", "Mohamed Farseen C S", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "required model not available"),
    @(58, 56, "21-12-2025", "vishnu", 8129349296, "29-12-2025", "SWAROOP S K", "Loss", "CUSTOMER INTERNAL ISSUES", "FAMILY DISAPPROVEL", "-", "Family problem Will be come"),
    @(59, 57, "21-12-2025", "thanveer", 9947874867, "07-01-2026", "Binu Augustine", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "."),
    @(60, 58, "21-12-2025", "tom", 8281692848, "03-01-2026", "Binu Augustine", "Loss", "ENQUIRY", "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "."),
    @(61, 59, "22-12-2025", "Madhav", 6282156345, "28-12-2025", "Althaf A", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", "black size 38 required"),
    @(62, 60, "22-12-2025", "Tony", 8075613454, "14-01-2026", "Althaf A", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", "customer required date product is not available"),
    @(63, 61, "22-12-2025", "noyal", 7034875397, "10-01-2026", "AKSHAYDAS.M", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", ".."),
    @(64, 62, "22-12-2025", "shahul", 8921029730, "18-01-2026", "AKSHAYDAS.M", "Loss", "PRODUCT", "PRODUCT NOT AVAILABLE", "-", ".."),
    @(65, 63, "23-12-2025", "Anandu", 9526656855, "08-01-2026", "Althaf A", "Loss", "SIZE NOT SUITABLE", "SIZE TOO SMALL", "-", "required 46 premium suit"),
    @(66, 64, "24-12-2025", "Athul", 9744852238, "07-02-2026", "Althaf A", "Loss", "PRODUCT", "REQUIRED MODEL NOT AVAILABLE", "-", "required design and color not available ( black indowestern with stone work )"),
    @(67, 65, "25-12-2025", "althaf", 9656287663, "11-01-2026", "ABIN KRISHNA S", "Loss", "ENQUIRY", "Enquiry for Relative/Friend", "-", "enquiry")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    $ws.Range("F$rowNum").Value = $r[6]
    $ws.Range("G$rowNum").Value = $r[7]
    $ws.Range("H$rowNum").Value = $r[8]
    $ws.Range("I$rowNum").Value = $r[9]
    $ws.Range("J$rowNum").Value = $r[10]
    if ($r[11] -ne $null) {
        $ws.Range("K$rowNum").Value = $r[11]
    } else {
        $ws.Range("K$rowNum").Value = ""
    }
    $ws.Range("A$rowNum").NumberFormat = "0"
    $ws.Range("D$rowNum").NumberFormat = "0"
}